$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Move the existing totals row (currently row 5) down to row 8 ---
$totalsLabel = $ws.Range("A5").Value()
$totalsFormula = $ws.Range("B5").Formula
$ws.Range("A5:B5").ClearContents()

# --- Grow the table (Tabell3) so it spans the two new activity rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C8"))

# --- New activity rows (5 and 6) ---
$ws.Range("A5").Value = 45307
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "Hämta alla aktiviteter"

$ws.Range("A6").Value = 45307
$ws.Range("A6").NumberFormat = $ws.Range("A4").NumberFormat
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "Spara aktivitet"

# --- Restore the totals row at its new location (row 8, leaving row 7 blank) ---
$ws.Range("A8").Value = $totalsLabel
$ws.Range("B8").Formula = $totalsFormula

# --- Update selection to reflect where the user left off ---
$ws.Range("A7").Select()
